# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml) for
# its slide master/design. The target revision swaps the theme colour
# scheme so the deck uses the stock "Office Theme" palette instead (the
# palette that used to live only on ppt/theme/theme1.xml, which is wired
# to the notes master).
#
# PowerPoint's Design > Colors gallery edits the active design's theme
# colour scheme in place (ppt/theme/theme2.xml here) without touching any
# relationships, which is exactly the effect we need.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Office Theme colour scheme values (RGB packed as &H00BBGGRR, matching
# the VBA ColorScheme.Colors(...).RGB convention):
$cs.Item(1).RGB  = 0          # Dark 1    - 000000
$cs.Item(2).RGB  = 16777215   # Light 1   - FFFFFF
$cs.Item(3).RGB  = 6968388    # Dark 2    - 44546A
$cs.Item(4).RGB  = 15132391   # Light 2   - E7E6E6
$cs.Item(5).RGB  = 13998939   # Accent 1  - 5B9BD5
$cs.Item(6).RGB  = 3243501    # Accent 2  - ED7D31
$cs.Item(7).RGB  = 10855845   # Accent 3  - A5A5A5
$cs.Item(8).RGB  = 49407      # Accent 4  - FFC000
$cs.Item(9).RGB  = 12874308   # Accent 5  - 4472C4
$cs.Item(10).RGB = 4697456    # Accent 6  - 70AD47
$cs.Item(11).RGB = 12673797   # Hyperlink - 0563C1
$cs.Item(12).RGB = 7491477    # Followed Hyperlink - 954F72
